# Auto-generated edit script: updates FFXIV Leve profit market-data values
# across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the
# scraped/refreshed market data, per the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 7975.5
$ws.Range("I19").Value = 10951
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 10951
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = -10776
$ws.Range("N19").Value = -5350
# Row 132
$ws.Range("H132").Value = 2569.4827
$ws.Range("I132").Value = 2621.3208
$ws.Range("K132").Value = 7863.9624
$ws.Range("M132").Value = -5333.9624
# Row 137
$ws.Range("H137").Value = 21740322
$ws.Range("I137").Value = 1019.74286
$ws.Range("J137").Value = 90910824
$ws.Range("K137").Value = 3059.22858
$ws.Range("L137").Value = 272732472
$ws.Range("M137").Value = -509.22858
$ws.Range("N137").Value = -272737572
# Row 138
$ws.Range("H138").Value = 2182.4087
$ws.Range("I138").Value = 1628.9434
$ws.Range("J138").Value = 2915.75
$ws.Range("K138").Value = 4886.8302
$ws.Range("L138").Value = 8747.25
$ws.Range("M138").Value = 253.1697999999997
$ws.Range("N138").Value = -19027.25
# Row 141
$ws.Range("H141").Value = 1178.4222
$ws.Range("I141").Value = 510.70587
$ws.Range("J141").Value = 3242.2727
$ws.Range("K141").Value = 1532.11761
$ws.Range("L141").Value = 9726.8181
$ws.Range("M141").Value = 3647.88239
$ws.Range("N141").Value = -20086.8181

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3297.13
$ws.Range("I32").Value = 3012.4675
$ws.Range("J32").Value = 4250.1304
$ws.Range("K32").Value = 3012.4675
$ws.Range("L32").Value = 4250.1304
$ws.Range("M32").Value = -2725.4675
$ws.Range("N32").Value = -4824.1304
# Row 74
$ws.Range("H74").Value = 10207803
$ws.Range("I74").Value = 13158626
$ws.Range("J74").Value = 14051
$ws.Range("K74").Value = 13158626
$ws.Range("L74").Value = 14051
$ws.Range("M74").Value = -13157752
$ws.Range("N74").Value = -15799
# Row 76
$ws.Range("H76").Value = 19429.334
$ws.Range("J76").Value = 19429.334
$ws.Range("L76").Value = 19429.334
$ws.Range("N76").Value = -20105.334
# Row 77
$ws.Range("H77").Value = 10207803
$ws.Range("I77").Value = 13158626
$ws.Range("J77").Value = 14051
$ws.Range("K77").Value = 65793130
$ws.Range("L77").Value = 70255
$ws.Range("M77").Value = -65788762
$ws.Range("N77").Value = -78991
# Row 79
$ws.Range("H79").Value = 19429.334
$ws.Range("J79").Value = 19429.334
$ws.Range("L79").Value = 19429.334
$ws.Range("N79").Value = -21769.334
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# Row 132
$ws.Range("H132").Value = 11254170
$ws.Range("I132").Value = 13547672
$ws.Range("K132").Value = 40643016
$ws.Range("M132").Value = -40640486
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 135
$ws.Range("H135").Value = 38667.06
$ws.Range("J135").Value = 38667.06
$ws.Range("L135").Value = 38667.06
$ws.Range("N135").Value = -48807.06
# Row 140
$ws.Range("H140").Value = 62863
$ws.Range("J140").Value = 62863
$ws.Range("L140").Value = 62863
$ws.Range("N140").Value = -73223

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5592937
$ws.Range("I31").Value = 1244.5172
$ws.Range("J31").Value = 46132708
$ws.Range("K31").Value = 1244.5172
$ws.Range("L31").Value = 46132708
$ws.Range("M31").Value = -949.5172
$ws.Range("N31").Value = -46133298
# Row 34
$ws.Range("H34").Value = 5592937
$ws.Range("I34").Value = 1244.5172
$ws.Range("J34").Value = 46132708
$ws.Range("K34").Value = 1244.5172
$ws.Range("L34").Value = 46132708
$ws.Range("M34").Value = -1042.5172
$ws.Range("N34").Value = -46133112
# Row 58
$ws.Range("H58").Value = 1186.2667
$ws.Range("I58").Value = 888.5
$ws.Range("J58").Value = 2202.1765
$ws.Range("K58").Value = 888.5
$ws.Range("L58").Value = 2202.1765
$ws.Range("M58").Value = -685.5
$ws.Range("N58").Value = -2608.1765
# Row 134
$ws.Range("H134").Value = 1239.1538
$ws.Range("I134").Value = 1417.1786
$ws.Range("J134").Value = 786
$ws.Range("K134").Value = 4251.5358
$ws.Range("L134").Value = 2358
$ws.Range("M134").Value = -1716.5358
$ws.Range("N134").Value = -7428
# Row 136
$ws.Range("H136").Value = 1186.2667
$ws.Range("I136").Value = 888.5
$ws.Range("J136").Value = 2202.1765
$ws.Range("K136").Value = 2665.5
$ws.Range("L136").Value = 6606.529500000001
$ws.Range("M136").Value = -115.5
$ws.Range("N136").Value = -11706.5295

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 2213.8333
$ws.Range("I34").Value = 199.875
$ws.Range("J34").Value = 3220.8125
$ws.Range("K34").Value = 599.625
$ws.Range("L34").Value = 9662.4375
$ws.Range("M34").Value = -515.625
$ws.Range("N34").Value = -9830.4375
# Row 113
$ws.Range("H113").Value = 495.58823
$ws.Range("I113").Value = 493.975
$ws.Range("J113").Value = 497.89285
$ws.Range("K113").Value = 1481.925
$ws.Range("L113").Value = 1493.67855
$ws.Range("M113").Value = 688.0749999999998
$ws.Range("N113").Value = -5833.678550000001
# Row 131
$ws.Range("H131").Value = 4669.143
$ws.Range("J131").Value = 3690.3333
$ws.Range("L131").Value = 11070.9999
$ws.Range("N131").Value = -21150.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 34
$ws.Range("H34").Value = 28000
$ws.Range("J34").Value = 28000
$ws.Range("L34").Value = 28000
$ws.Range("N34").Value = -28536
# Row 76
$ws.Range("H76").Value = 28000
$ws.Range("J76").Value = 28000
$ws.Range("L76").Value = 28000
$ws.Range("N76").Value = -28630
# Row 79
$ws.Range("H79").Value = 28000
$ws.Range("J79").Value = 28000
$ws.Range("L79").Value = 28000
$ws.Range("N79").Value = -30184
# Row 80
$ws.Range("H80").Value = 36358.594
$ws.Range("I80").Value = 2268.182
$ws.Range("J80").Value = 54215.477
$ws.Range("K80").Value = 2268.182
$ws.Range("L80").Value = 54215.477
$ws.Range("M80").Value = -1270.182
$ws.Range("N80").Value = -56211.477
# Row 83
$ws.Range("H83").Value = 36358.594
$ws.Range("I83").Value = 2268.182
$ws.Range("J83").Value = 54215.477
$ws.Range("K83").Value = 11340.91
$ws.Range("L83").Value = 271077.385
$ws.Range("M83").Value = -6348.91
$ws.Range("N83").Value = -281061.385
# Row 97
$ws.Range("H97").Value = 1989.25
$ws.Range("I97").Value = 2016.3636
$ws.Range("J97").Value = 1691
$ws.Range("K97").Value = 2016.3636
$ws.Range("L97").Value = 1691
$ws.Range("M97").Value = -1520.3636
$ws.Range("N97").Value = -2683
# Row 132
$ws.Range("H132").Value = 40002424
$ws.Range("I132").Value = 62502320
$ws.Range("J132").Value = 2610.889
$ws.Range("K132").Value = 187506960
$ws.Range("L132").Value = 7832.667
$ws.Range("M132").Value = -187504430
$ws.Range("N132").Value = -12892.667

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2038.6578
$ws.Range("I16").Value = 1856.3928
$ws.Range("K16").Value = 1856.3928
$ws.Range("M16").Value = -1686.3928
# Row 40
$ws.Range("H40").Value = 1501
$ws.Range("I40").Value = 1501
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1501
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1365
$ws.Range("N40").ClearContents()
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 136
$ws.Range("H136").Value = 948.8125
$ws.Range("I136").Value = 530.6070999999999
$ws.Range("J136").Value = 3876.25
$ws.Range("K136").Value = 1591.8213
$ws.Range("L136").Value = 11628.75
$ws.Range("M136").Value = 958.1787000000002
$ws.Range("N136").Value = -16728.75

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 7693274
$ws.Range("I132").Value = 8718626
$ws.Range("J132").Value = 3133
$ws.Range("K132").Value = 26155878
$ws.Range("L132").Value = 9399
$ws.Range("M132").Value = -26153348
$ws.Range("N132").Value = -14459
